$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column D
$ws.Range("D1").Value = "peak2"

# D2 standalone formula
$ws.Range("D2").Formula = "=B2*1.13"

# D3:D53 shared formula block
$ws.Range("D3:D53").Formula = "=B3*1.13"

# Column B width
$ws.Columns("B").ColumnWidth = 17.1640625

# View settings
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("D2:D53").Select()
